$o = New-Object -ComObject "Word.Basic"
Write-Output $o
